$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1 / table3) - add rows 4 & 5
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "c70687d6-445d-4c97-85ed-ea89b6e51f0c.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-09-02 14:50:41"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add($wsOverview.Cells(4,2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/HEAD/e2e/c70687d6-445d-4c97-85ed-ea89b6e51f0c.md", "", "", 'e2e\c70687d6-445d-4c97-85ed-ea89b6e51f0c.md')

$wsOverview.Range("A5").Value = "d3f2995f-e8d3-4e47-909f-09b781ef0592.md"
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-09-02 14:50:41"
$wsOverview.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add($wsOverview.Cells(5,2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/HEAD/e2e/d3f2995f-e8d3-4e47-909f-09b781ef0592.md", "", "", 'e2e\d3f2995f-e8d3-4e47-909f-09b781ef0592.md')

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G5"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2 / table1) - add rows 4 & 5
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'False"
$wsZh.Range("G4").Value = "c70687d6-445d-4c97-85ed-ea89b6e51f0c.dec5f9343355d0db5eacbce61db3ef90e90315ee.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-09-02 14:50:36"
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K4").Value = "0001-01-01 00:00:00"
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("O4").Value = "'False"
$wsZh.Hyperlinks.Add($wsZh.Cells(4,1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/HEAD/e2e/c70687d6-445d-4c97-85ed-ea89b6e51f0c.md", "", "", "c70687d6-445d-4c97-85ed-ea89b6e51f0c.md")

$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "'False"
$wsZh.Range("G5").Value = "d3f2995f-e8d3-4e47-909f-09b781ef0592.dc128ba0ed4c5c896cc6873ff954d6ea76266e1a.zh-cn.xlf"
$wsZh.Range("H5").Value = "2016-09-02 14:50:36"
$wsZh.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K5").Value = "0001-01-01 00:00:00"
$wsZh.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("M5").Value = "'True"
$wsZh.Range("O5").Value = "'False"
$wsZh.Hyperlinks.Add($wsZh.Cells(5,1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/HEAD/e2e/d3f2995f-e8d3-4e47-909f-09b781ef0592.md", "", "", "d3f2995f-e8d3-4e47-909f-09b781ef0592.md")

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P5"))

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3 / table2) - add rows 4 & 5
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'False"
$wsDe.Range("G4").Value = "c70687d6-445d-4c97-85ed-ea89b6e51f0c.dec5f9343355d0db5eacbce61db3ef90e90315ee.de-de.xlf"
$wsDe.Range("H4").Value = "2016-09-02 14:50:41"
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K4").Value = "0001-01-01 00:00:00"
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("O4").Value = "'False"
$wsDe.Hyperlinks.Add($wsDe.Cells(4,1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/HEAD/e2e/c70687d6-445d-4c97-85ed-ea89b6e51f0c.md", "", "", "c70687d6-445d-4c97-85ed-ea89b6e51f0c.md")

$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "'False"
$wsDe.Range("G5").Value = "d3f2995f-e8d3-4e47-909f-09b781ef0592.dc128ba0ed4c5c896cc6873ff954d6ea76266e1a.de-de.xlf"
$wsDe.Range("H5").Value = "2016-09-02 14:50:41"
$wsDe.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDe.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("M5").Value = "'True"
$wsDe.Range("O5").Value = "'False"
$wsDe.Hyperlinks.Add($wsDe.Cells(5,1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/HEAD/e2e/d3f2995f-e8d3-4e47-909f-09b781ef0592.md", "", "", "d3f2995f-e8d3-4e47-909f-09b781ef0592.md")

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P5"))
